$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from the existing H1 header
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

$data = @{
    2  = @(7, 7)
    3  = @(8, 8)
    4  = @(7, 7)
    5  = @(7, 7)
    6  = @(9, 9)
    7  = @(7, 7)
    8  = @(8, 8)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(7, 7)
    16 = @(7, 7)
    17 = @(7, 7)
    18 = @(8, 8)
    19 = @(7, 8)
    20 = @(8, 8)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(7, 7)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(7, 7)
    27 = @(6, 6)
    28 = @(8, 8)
    29 = @(7, 7)
    30 = @(7, 7)
    31 = @(8, 8)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(8, 8)
    35 = @(7, 7)
    36 = @(5, 6)
    37 = @(7, 7)
    38 = @(9, 9)
    39 = @(6, 6)
    40 = @(9, 9)
    41 = @(7, 7)
    42 = @(7, 7)
    43 = @(8, 8)
    44 = @(8, 8)
    45 = @(6, 7)
    46 = @(7, 7)
    47 = @(6, 6)
    48 = @(8, 8)
    49 = @(7, 7)
    50 = @(9, 9)
    51 = @(8, 8)
    52 = @(5, 5)
    53 = @(8, 8)
    54 = @(9, 9)
    55 = @(7, 7)
    56 = @(8, 8)
    57 = @(8, 8)
    58 = @(8, 8)
    59 = @(9, 9)
    60 = @(8, 8)
    61 = @(6, 6)
    62 = @(7, 7)
    63 = @(6, 6)
    64 = @(6, 6)
    65 = @(7, 7)
    66 = @(4, 4)
    67 = @(3, 3)
    68 = @(4, 4)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
